$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 21
$ws.Range("I2").Value = 51
$ws.Range("J2").Value = 215
$ws.Range("K2").Value = 1
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 35
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 5
$ws.Range("S2").Value = 22
$ws.Range("T2").Value = 41
$ws.Range("V2").Value = 393
$ws.Range("X2").Value = 339
$ws.Range("Z2").Value = 5
$ws.Range("AA2").Value = 2

$wb.Save()
